$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new row 4 data (hours worked, week, hours needed, notes)
$ws.Range("A4").Value = 7
$ws.Range("C4").Value = 48
$ws.Range("D4").Value = "Worked on V2 of TTT and Stopwatch. Finished more of the quizzes. 1 hour on Monday 4 on Thursday 2 on Sunday"

# Widen column D to fit the new, longer note text
# (engine quantizes ColumnWidth to an MDW-7 pixel grid, same as real Excel COM;
# 79.3 is the input that lands closest to the target stored width of 80.21875)
$ws.Columns.Item(4).ColumnWidth = 79.3

# Update the active cell selection to D4
$ws.Range("D4").Select()
